# Auto-generated PowerShell/Excel-COM script
# Updates the '想去人数' (interested-count) column F
# across sheets 展览 / 演出 / 全部类型 to match the refreshed data snapshot.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 171
$ws1.Range("F4").Value = 398
$ws1.Range("F5").Value = 945
$ws1.Range("F6").Value = 5005
$ws1.Range("F7").Value = 396
$ws1.Range("F8").Value = 572
$ws1.Range("F9").Value = 871
$ws1.Range("F11").Value = 64
$ws1.Range("F12").Value = 17
$ws1.Range("F14").Value = 13
$ws1.Range("F16").Value = 1639
$ws1.Range("F17").Value = 1430
$ws1.Range("F18").Value = 723
$ws1.Range("F20").Value = 177
$ws1.Range("F21").Value = 265
$ws1.Range("F22").Value = 473
$ws1.Range("F23").Value = 117
$ws1.Range("F24").Value = 1039
$ws1.Range("F27").Value = 1929
$ws1.Range("F28").Value = 152
$ws1.Range("F29").Value = 83
$ws1.Range("F30").Value = 14
$ws1.Range("F31").Value = 206
$ws1.Range("F33").Value = 38
$ws1.Range("F36").Value = 257
$ws1.Range("F37").Value = 565
$ws1.Range("F38").Value = 72
$ws1.Range("F39").Value = 27
$ws1.Range("F40").Value = 27
$ws1.Range("F41").Value = 50

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 136
$ws2.Range("F9").Value = 10

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 171
$ws4.Range("F5").Value = 398
$ws4.Range("F6").Value = 945
$ws4.Range("F8").Value = 5005
$ws4.Range("F9").Value = 396
$ws4.Range("F10").Value = 572
$ws4.Range("F12").Value = 136
$ws4.Range("F13").Value = 871
$ws4.Range("F17").Value = 64
$ws4.Range("F18").Value = 17
$ws4.Range("F20").Value = 13
$ws4.Range("F23").Value = 1639
$ws4.Range("F24").Value = 1430
$ws4.Range("F25").Value = 723
$ws4.Range("F27").Value = 177
$ws4.Range("F28").Value = 265
$ws4.Range("F30").Value = 473
$ws4.Range("F31").Value = 117
$ws4.Range("F32").Value = 1039
$ws4.Range("F34").Value = 1929
$ws4.Range("F35").Value = 152
$ws4.Range("F36").Value = 83
$ws4.Range("F37").Value = 14
$ws4.Range("F38").Value = 206
$ws4.Range("F40").Value = 38
$ws4.Range("F42").Value = 257
$ws4.Range("F43").Value = 565
$ws4.Range("F44").Value = 72
$ws4.Range("F45").Value = 27
$ws4.Range("F46").Value = 27
$ws4.Range("F47").Value = 50
$ws4.Range("F48").Value = 10

